$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.234.81'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.71%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.656.40'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.00%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.63%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '219.24'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.77%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5238'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.34%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.005'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.60%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2656'

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.31%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.70'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.81%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.07%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.565'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.29%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.650.14'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.52%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.885.74'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.86%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5648'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.09%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8104'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.27%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.46'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.50%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.234.85'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.79%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.004'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.60%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.728'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.23%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '193.55'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.61%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.29'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.63%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.030'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.70%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.62%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.93'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.58%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.271'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.00%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.07'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.29%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.506'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.44%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05606'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -5.11%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.278'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.22%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.501'

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.381'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.05%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.593'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -2.24%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.804'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.69%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9441'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.06%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.19%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5763'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.33%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.78%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.944'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.61%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.589'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.47%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8476'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.42%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.039.13'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -3.63%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.95%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.795.76'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.88%  '

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.02%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.52%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05317'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.86%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.9982'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.64%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4355'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.07%  '
